# Update the "Förändrad" (Changed) date column (C) for rows 2-27
# from serial date 45266 (2023-12-06) to 45267 (2023-12-07).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

for ($row = 2; $row -le 27; $row++) {
    $cell = $ws.Cells.Item($row, 3)
    if ($cell.Value2 -eq 45266) {
        $cell.Value = 45267
    }
}
